$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H21").Value = 31379.75
$ws.Range("I21").Value = 80019
$ws.Range("J21").Value = 15166.667
$ws.Range("K21").Value = 80019
$ws.Range("L21").Value = 15166.667
$ws.Range("M21").Value = -79551
$ws.Range("N21").Value = -16102.667

$ws.Range("H23").Value = 31379.75
$ws.Range("I23").Value = 80019
$ws.Range("J23").Value = 15166.667
$ws.Range("K23").Value = 80019
$ws.Range("L23").Value = 15166.667
$ws.Range("M23").Value = -79785
$ws.Range("N23").Value = -15634.667

$ws.Range("H29").Value = 3433.3333
$ws.Range("I29").Value = 300
$ws.Range("J29").Value = 5000
$ws.Range("K29").Value = 900
$ws.Range("L29").Value = 15000
$ws.Range("M29").Value = -619
$ws.Range("N29").Value = -15562

$ws.Range("H38").Value = 1195.4706
$ws.Range("I38").Value = 194.22223
$ws.Range("J38").Value = 2321.875
$ws.Range("K38").Value = 582.66669
$ws.Range("L38").Value = 6965.625
$ws.Range("M38").Value = -210.66669
$ws.Range("N38").Value = -7709.625

$ws.Range("H58").Value = 1360.75
$ws.Range("J58").Value = 2940
$ws.Range("L58").Value = 8820
$ws.Range("N58").Value = -9120

$ws.Range("H87").Value = 25063.6
$ws.Range("J87").Value = 25063.6
$ws.Range("L87").Value = 25063.6
$ws.Range("N87").Value = -27559.6

$ws.Range("H90").Value = 25063.6
$ws.Range("J90").Value = 25063.6
$ws.Range("L90").Value = 75190.79999999999
$ws.Range("N90").Value = -87670.79999999999

$ws.Range("H132").Value = 4257491.5
$ws.Range("I132").Value = 4446475.5
$ws.Range("K132").Value = 13339426.5
$ws.Range("M132").Value = -13336896.5

$ws.Range("H135").Value = 961.1539
$ws.Range("I135").Value = 913.2857
$ws.Range("J135").Value = 1162.2
$ws.Range("K135").Value = 8219.5713
$ws.Range("L135").Value = 10459.8
$ws.Range("M135").Value = -5684.5713
$ws.Range("N135").Value = -15529.8


$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4361.989
$ws.Range("I32").Value = 3502.1646
$ws.Range("J32").Value = 10537.091
$ws.Range("K32").Value = 3502.1646
$ws.Range("L32").Value = 10537.091
$ws.Range("M32").Value = -3215.1646
$ws.Range("N32").Value = -11111.091

$ws.Range("H61").Value = 3722.0557
$ws.Range("I61").Value = 1624.625
$ws.Range("K61").Value = 1624.625
$ws.Range("M61").Value = -1412.625

$ws.Range("H136").Value = 3722.0557
$ws.Range("I136").Value = 1624.625
$ws.Range("K136").Value = 4873.875
$ws.Range("M136").Value = -2323.875


$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 1094.6923
$ws.Range("I94").Value = 863.1
$ws.Range("J94").Value = 1866.6666
$ws.Range("K94").Value = 863.1
$ws.Range("L94").Value = 1866.6666
$ws.Range("M94").Value = -412.1
$ws.Range("N94").Value = -2768.6666

$ws.Range("H105").Value = 2515.9285
$ws.Range("I105").Value = 2425
$ws.Range("K105").Value = 2425
$ws.Range("M105").Value = -678

$ws.Range("H123").Value = 0
$ws.Range("J123").Value = 0
$ws.Range("N123").ClearContents()

$ws.Range("H134").Value = 3052.3044
$ws.Range("I134").Value = 2544.5625
$ws.Range("K134").Value = 7633.6875
$ws.Range("M134").Value = -5098.6875


$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 12502738
$ws.Range("I58").Value = 1450.25
$ws.Range("J58").Value = 41672410
$ws.Range("K58").Value = 1450.25
$ws.Range("L58").Value = 41672410
$ws.Range("M58").Value = -1247.25
$ws.Range("N58").Value = -41672816

$ws.Range("H123").Value = 31951.111
$ws.Range("J123").Value = 31951.111
$ws.Range("L123").Value = 31951.111
$ws.Range("N123").Value = -41751.111

$ws.Range("H132").Value = 2506.6584
$ws.Range("I132").Value = 1837.8214
$ws.Range("K132").Value = 5513.4642
$ws.Range("M132").Value = -2983.4642

$ws.Range("H134").Value = 2145.889
$ws.Range("I134").Value = 1087.4286
$ws.Range("K134").Value = 3262.2858
$ws.Range("M134").Value = -727.2857999999997

$ws.Range("H136").Value = 12502738
$ws.Range("I136").Value = 1450.25
$ws.Range("J136").Value = 41672410
$ws.Range("K136").Value = 4350.75
$ws.Range("L136").Value = 125017230
$ws.Range("M136").Value = -1800.75
$ws.Range("N136").Value = -125022330


$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H74").Value = 10460
$ws.Range("J74").Value = 11187.692
$ws.Range("L74").Value = 33563.076
$ws.Range("N74").Value = -35685.076

$ws.Range("H77").Value = 10460
$ws.Range("J77").Value = 11187.692
$ws.Range("L77").Value = 100689.228
$ws.Range("N77").Value = -111297.228

$ws.Range("H81").Value = 10183.333
$ws.Range("I81").Value = 300
$ws.Range("J81").Value = 15125
$ws.Range("K81").Value = 900
$ws.Range("L81").Value = 45375
$ws.Range("M81").Value = 223
$ws.Range("N81").Value = -47621

$ws.Range("H84").Value = 10183.333
$ws.Range("I84").Value = 300
$ws.Range("J84").Value = 15125
$ws.Range("K84").Value = 2700
$ws.Range("L84").Value = 136125
$ws.Range("M84").Value = 2916
$ws.Range("N84").Value = -147357

$ws.Range("H131").Value = 1053.5369
$ws.Range("I131").Value = 1270.1333
$ws.Range("J131").Value = 1012.925
$ws.Range("K131").Value = 3810.3999
$ws.Range("L131").Value = 3038.775
$ws.Range("M131").Value = 1229.6001
$ws.Range("N131").Value = -13118.775

$ws.Range("H141").Value = 4833.3335
$ws.Range("I141").Value = 3666.6667
$ws.Range("J141").Value = 6000
$ws.Range("K141").Value = 11000.0001
$ws.Range("L141").Value = 18000
$ws.Range("M141").Value = -5820.000100000001
$ws.Range("N141").Value = -28360


$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 4848.5713
$ws.Range("I70").Value = 4908.3335
$ws.Range("J70").Value = 4490
$ws.Range("K70").Value = 4908.3335
$ws.Range("L70").Value = 4490
$ws.Range("M70").Value = -4638.3335
$ws.Range("N70").Value = -5030

$ws.Range("H73").Value = 4848.5713
$ws.Range("I73").Value = 4908.3335
$ws.Range("J73").Value = 4490
$ws.Range("K73").Value = 4908.3335
$ws.Range("L73").Value = 4490
$ws.Range("M73").Value = -3972.3335
$ws.Range("N73").Value = -6362

$ws.Range("H80").Value = 4499.625
$ws.Range("I80").Value = 4666.5
$ws.Range("J80").Value = 3999
$ws.Range("K80").Value = 4666.5
$ws.Range("L80").Value = 3999
$ws.Range("M80").Value = -3668.5
$ws.Range("N80").Value = -5995

$ws.Range("H83").Value = 4499.625
$ws.Range("I83").Value = 4666.5
$ws.Range("J83").Value = 3999
$ws.Range("K83").Value = 23332.5
$ws.Range("L83").Value = 19995
$ws.Range("M83").Value = -18340.5
$ws.Range("N83").Value = -29979

$ws.Range("H102").Value = 114777.664
$ws.Range("I102").Value = 2833.3333
$ws.Range("J102").Value = 338666.34
$ws.Range("K102").Value = 2833.3333
$ws.Range("L102").Value = 338666.34
$ws.Range("M102").Value = -1211.3333
$ws.Range("N102").Value = -341910.34

$ws.Range("H122").Value = 3927.6924
$ws.Range("I122").Value = 2350
$ws.Range("J122").Value = 4628.8887
$ws.Range("K122").Value = 7050
$ws.Range("L122").Value = 13886.6661
$ws.Range("M122").Value = -4600
$ws.Range("N122").Value = -18786.6661


$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 55556396
$ws.Range("I22").Value = 100000280
$ws.Range("J22").Value = 1543.625
$ws.Range("K22").Value = 100000280
$ws.Range("L22").Value = 1543.625
$ws.Range("M22").Value = -99999985
$ws.Range("N22").Value = -2133.625

$ws.Range("H27").Value = 55556396
$ws.Range("I27").Value = 100000280
$ws.Range("J27").Value = 1543.625
$ws.Range("K27").Value = 100000280
$ws.Range("L27").Value = 1543.625
$ws.Range("M27").Value = -100000173
$ws.Range("N27").Value = -1757.625

$ws.Range("H46").Value = 2114.9167
$ws.Range("I46").Value = 498.33334
$ws.Range("J46").Value = 3731.5
$ws.Range("K46").Value = 498.33334
$ws.Range("L46").Value = 3731.5
$ws.Range("M46").Value = -310.33334
$ws.Range("N46").Value = -4107.5

$ws.Range("H55").Value = 1124.6428
$ws.Range("I55").Value = 325.16666
$ws.Range("J55").Value = 1724.25
$ws.Range("K55").Value = 325.16666
$ws.Range("L55").Value = 1724.25
$ws.Range("M55").Value = -152.16666
$ws.Range("N55").Value = -2070.25

$ws.Range("H68").Value = 1421.0526
$ws.Range("I68").Value = 1000
$ws.Range("K68").Value = 1000
$ws.Range("M68").Value = -251

$ws.Range("H71").Value = 1421.0526
$ws.Range("I71").Value = 1000
$ws.Range("K71").Value = 5000
$ws.Range("M71").Value = -1256

$ws.Range("H93").Value = 4091.111
$ws.Range("I93").Value = 1740.3334
$ws.Range("J93").Value = 5266.5
$ws.Range("K93").Value = 1740.3334
$ws.Range("L93").Value = 5266.5
$ws.Range("M93").Value = -492.3334
$ws.Range("N93").Value = -7762.5

$ws.Range("H123").Value = 23400
$ws.Range("J123").Value = 23400
$ws.Range("L123").Value = 23400
$ws.Range("N123").Value = -33200


$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 2577.6155
$ws.Range("I107").Value = 2650.375
$ws.Range("J107").Value = 2461.2
$ws.Range("K107").Value = 7951.125
$ws.Range("L107").Value = 7383.599999999999
$ws.Range("M107").Value = -6031.125
$ws.Range("N107").Value = -11223.6

$ws.Range("H132").Value = 167408.25
$ws.Range("I132").Value = 251506.25
$ws.Range("J132").Value = 7221.5713
$ws.Range("K132").Value = 754518.75
$ws.Range("L132").Value = 21664.7139
$ws.Range("M132").Value = -751988.75
$ws.Range("N132").Value = -26724.7139
